# Auto-generated edit script applying cell-value changes per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 16667879
$ws.Range("I6").Value = 47619144
$ws.Range("K6").Value = 142857432
$ws.Range("M6").Value = -142857320
$ws.Range("H100").Value = 166262.17
$ws.Range("I100").Value = 188168.25
$ws.Range("K100").Value = 188168.25
$ws.Range("M100").Value = -187627.25
$ws.Range("H106").Value = 6502394
$ws.Range("I106").Value = 8822808
$ws.Range("K106").Value = 8822808
$ws.Range("M106").Value = -8822177
$ws.Range("H125").Value = 6429.3
$ws.Range("I125").Value = 7042.4287
$ws.Range("K125").Value = 63381.85830000001
$ws.Range("M125").Value = -60921.85830000001
$ws.Range("H132").Value = 4768628
$ws.Range("I132").Value = 6954.9443
$ws.Range("K132").Value = 20864.8329
$ws.Range("M132").Value = -18334.8329
$ws.Range("H135").Value = 5991.8486
$ws.Range("I135").Value = 8579.866
$ws.Range("J135").Value = 3835.1667
$ws.Range("K135").Value = 77218.79399999999
$ws.Range("L135").Value = 34516.5003
$ws.Range("M135").Value = -74683.79399999999
$ws.Range("N135").Value = -39586.5003
$ws.Range("H138").Value = 367938.03
$ws.Range("I138").Value = 659206.9
$ws.Range("K138").Value = 1977620.7
$ws.Range("M138").Value = -1972480.7
$ws.Range("H141").Value = 6017.143
$ws.Range("I141").Value = 5898.095
$ws.Range("K141").Value = 17694.285
$ws.Range("M141").Value = -12514.285

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5000
$ws.Range("I26").Value = 5000
$ws.Range("K26").Value = 5000
$ws.Range("M26").Value = -4670
$ws.Range("H32").Value = 4509.983
$ws.Range("I32").Value = 4210.5273
$ws.Range("K32").Value = 4210.5273
$ws.Range("M32").Value = -3923.5273
$ws.Range("H61").Value = 15465.134
$ws.Range("I61").Value = 26141.428
$ws.Range("K61").Value = 26141.428
$ws.Range("M61").Value = -25929.428
$ws.Range("H88").Value = 1618.1428
$ws.Range("I88").Value = 960.6
$ws.Range("J88").Value = 1983.4445
$ws.Range("K88").Value = 960.6
$ws.Range("L88").Value = 1983.4445
$ws.Range("M88").Value = -554.6
$ws.Range("N88").Value = -2795.4445
$ws.Range("H91").Value = 1618.1428
$ws.Range("I91").Value = 960.6
$ws.Range("J91").Value = 1983.4445
$ws.Range("K91").Value = 960.6
$ws.Range("L91").Value = 1983.4445
$ws.Range("M91").Value = 443.4
$ws.Range("N91").Value = -4791.4445
$ws.Range("H119").Value = 23000
$ws.Range("J119").Value = 23000
$ws.Range("L119").Value = 23000
$ws.Range("N119").Value = -32676
$ws.Range("H132").Value = 6007.5
$ws.Range("I132").Value = 5515
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 16545
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -14015
$ws.Range("N132").Value = -24560
$ws.Range("H136").Value = 15465.134
$ws.Range("I136").Value = 26141.428
$ws.Range("K136").Value = 78424.284
$ws.Range("M136").Value = -75874.284

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3458
$ws.Range("I86").Value = 4643.478
$ws.Range("J86").Value = 1510.4286
$ws.Range("K86").Value = 4643.478
$ws.Range("L86").Value = 1510.4286
$ws.Range("M86").Value = -3520.478
$ws.Range("N86").Value = -3756.4286
$ws.Range("H89").Value = 3458
$ws.Range("I89").Value = 4643.478
$ws.Range("J89").Value = 1510.4286
$ws.Range("K89").Value = 23217.39
$ws.Range("L89").Value = 7552.143
$ws.Range("M89").Value = -17601.39
$ws.Range("N89").Value = -18784.143
$ws.Range("H96").Value = 9570.799999999999
$ws.Range("I96").Value = 9570.799999999999
$ws.Range("K96").Value = 9570.799999999999
$ws.Range("M96").Value = -6824.799999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 20456.8
$ws.Range("I7").Value = 40323.6
$ws.Range("J7").Value = 590
$ws.Range("K7").Value = 40323.6
$ws.Range("L7").Value = 590
$ws.Range("M7").Value = -40210.6
$ws.Range("N7").Value = -816
$ws.Range("H31").Value = 11263.071
$ws.Range("I31").Value = 12789.363
$ws.Range("K31").Value = 12789.363
$ws.Range("M31").Value = -12494.363
$ws.Range("H34").Value = 11263.071
$ws.Range("I34").Value = 12789.363
$ws.Range("K34").Value = 12789.363
$ws.Range("M34").Value = -12587.363
$ws.Range("H122").Value = 11677.417
$ws.Range("I122").Value = 16530.875
$ws.Range("K122").Value = 49592.625
$ws.Range("M122").Value = -47142.625
$ws.Range("H134").Value = 14585.125
$ws.Range("I134").Value = 18447
$ws.Range("K134").Value = 55341
$ws.Range("M134").Value = -52806
$ws.Range("H141").Value = 338345.16
$ws.Range("J141").Value = 411530.2
$ws.Range("L141").Value = 411530.2
$ws.Range("N141").Value = -421890.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39179880
$ws.Range("I4").Value = 41080864
$ws.Range("K4").Value = 123242592
$ws.Range("M4").Value = -123242480
$ws.Range("H12").Value = 29.076923
$ws.Range("I12").Value = 101.333336
$ws.Range("J12").Value = 7.4
$ws.Range("K12").Value = 304.000008
$ws.Range("L12").Value = 22.2
$ws.Range("M12").Value = -131.000008
$ws.Range("N12").Value = -368.2
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H98").Value = 2065.0833
$ws.Range("I98").Value = 2173.875
$ws.Range("J98").Value = 1847.5
$ws.Range("K98").Value = 6521.625
$ws.Range("L98").Value = 5542.5
$ws.Range("M98").Value = -5023.625
$ws.Range("N98").Value = -8538.5
$ws.Range("H126").Value = 14581.8
$ws.Range("J126").Value = 27555
$ws.Range("L126").Value = 82665
$ws.Range("N126").Value = -92545
$ws.Range("H129").Value = 23811994
$ws.Range("I129").Value = 860
$ws.Range("J129").Value = 37040400
$ws.Range("K129").Value = 2580
$ws.Range("L129").Value = 111121200
$ws.Range("M129").Value = 2420
$ws.Range("N129").Value = -111131200

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 29500
$ws.Range("H132").Value = 11824.75
$ws.Range("I132").Value = 13014
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 39042
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -36512
$ws.Range("N132").Value = -15560

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 22350.902
$ws.Range("I40").Value = 27049.611
$ws.Range("J40").Value = 15845
$ws.Range("K40").Value = 27049.611
$ws.Range("L40").Value = 15845
$ws.Range("M40").Value = -26913.611
$ws.Range("N40").Value = -16117
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H96").Value = 44999.5
$ws.Range("J96").Value = 44999.5
$ws.Range("L96").Value = 44999.5
$ws.Range("N96").Value = -50491.5
$ws.Range("H122").Value = 6850.65
$ws.Range("I122").Value = 7335
$ws.Range("J122").Value = 6258.6665
$ws.Range("K122").Value = 22005
$ws.Range("L122").Value = 18775.9995
$ws.Range("M122").Value = -19555
$ws.Range("N122").Value = -23675.9995
$ws.Range("H132").Value = 624496.5600000001
$ws.Range("I132").Value = 748146.1
$ws.Range("J132").Value = 6248.75
$ws.Range("K132").Value = 2244438.3
$ws.Range("L132").Value = 18746.25
$ws.Range("M132").Value = -2241908.3
$ws.Range("N132").Value = -23806.25
$ws.Range("H136").Value = 5645.778
$ws.Range("I136").Value = 4802
$ws.Range("K136").Value = 14406
$ws.Range("M136").Value = -11856

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 13984.5
$ws.Range("I58").Value = 13984.5
$ws.Range("K58").Value = 13984.5
$ws.Range("M58").Value = -13676.5
$ws.Range("H107").Value = 21207.223
$ws.Range("I107").Value = 5625.357
$ws.Range("J107").Value = 75743.75
$ws.Range("K107").Value = 16876.071
$ws.Range("L107").Value = 227231.25
$ws.Range("M107").Value = -14956.071
$ws.Range("N107").Value = -231071.25
$ws.Range("H126").Value = 36038.582
$ws.Range("I126").Value = 59053
$ws.Range("K126").Value = 177159
$ws.Range("M126").Value = -174689
$ws.Range("H132").Value = 10127.923
$ws.Range("I132").Value = 11586.827
$ws.Range("J132").Value = 5897.1
$ws.Range("K132").Value = 34760.481
$ws.Range("L132").Value = 17691.3
$ws.Range("M132").Value = -32230.481
$ws.Range("N132").Value = -22751.3
$ws.Range("H139").Value = 102749.75
$ws.Range("J139").Value = 66999.664
$ws.Range("L139").Value = 66999.664
$ws.Range("N139").Value = -77279.664
